$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.861.95"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.598.60"
$ws.Range("E3").Value = "  -2.11%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.13"
$ws.Range("E5").Value = "  -2.31%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.477"
$ws.Range("E7").Value = "  -5.48%  "

$ws.Range("E8").Value = "  -2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0611"
$ws.Range("E9").Value = "  -2.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.82"
$ws.Range("E10").Value = "  -3.71%  "

$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.14"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.604.17"
$ws.Range("E13").Value = "  -1.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  -3.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.508"
$ws.Range("E15").Value = "  -4.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.845.45"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.63"
$ws.Range("E17").Value = "  -1.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0714"
$ws.Range("E18").Value = "  -4.05%  "

$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.10"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.17"
$ws.Range("E21").Value = "  -1.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.28"
$ws.Range("E22").Value = "  -2.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  -3.26%  "

$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.61"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("E26").Value = "  -3.84%  "

$ws.Range("E27").Value = "  -3.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.50"
$ws.Range("E28").Value = "  -3.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.90"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("E30").Value = "  -2.64%  "

$ws.Range("E31").Value = "  -3.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.07"
$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("E33").Value = "  -4.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.40"
$ws.Range("E34").Value = "  -1.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.46"
$ws.Range("E35").Value = "  -2.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.107.28"
$ws.Range("E36").Value = "  -2.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -3.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.796"
$ws.Range("E38").Value = "  -8.64%  "

$ws.Range("E39").Value = "  -2.97%  "

$ws.Range("E40").Value = "  -5.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "95.55"
$ws.Range("E41").Value = "  -3.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.734.68"
$ws.Range("E42").Value = "  -2.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$ws.Range("E43").Value = "  -3.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.737"
$ws.Range("E44").Value = "  -5.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0113"
$ws.Range("E45").Value = "  -1.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "52.96"
$ws.Range("E46").Value = "  -3.89%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0510"
$ws.Range("E47").Value = "  -3.21%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.45"
$ws.Range("E48").Value = "  -1.84%  "

$ws.Range("E49").Value = "  -1.04%  "

$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.29"
$ws.Range("E51").Value = "  -3.39%  "

